$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - UNH
$ws.Range("D2").Value = 324.54
$ws.Range("E2").Value = 48
$ws.Range("F2").Value = 1.72
$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 43
$ws.Range("K2").Value = 57.8
$ws.Range("N2").Value = 66.04328690552585

# Row 3 - PRU
$ws.Range("D3").Value = 107.85
$ws.Range("E3").Value = 59.9
$ws.Range("F3").Value = 0.62
$ws.Range("K3").Value = 56.2
$ws.Range("N3").Value = 66.04328690552585

# Row 4 - MET
$ws.Range("D4").Value = 76.95
$ws.Range("E4").Value = 54.3
$ws.Range("F4").Value = 2.45
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 36
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 54.8
$ws.Range("N4").Value = 66.04328690552585

# Row 5 - AIG
$ws.Range("D5").Value = 76.27
$ws.Range("E5").Value = 43.5
$ws.Range("F5").Value = 1.44
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 53
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 53
$ws.Range("N5").Value = 66.04328690552585
